$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "gender" column is dropped. Rather than a structural delete, the data was
# re-entered shifted one column to the left (C<-D, D<-E, E<-F, F<-G), with some
# values/headers renamed, and the now-unused column G cleared out.

# Headers first: D1 (studyStage, was citizenshipStatus) and E1 (courseCode, was course).
$ws.Range("D1").Value = "studyStage"
$ws.Range("E1").Value = "courseCode"

# Remaining headers reuse existing strings.
$ws.Range("C1").Value = "citizenshipStatus"
$ws.Range("F1").Value = "pemGroup"

# Row values, column by column (C, D, E, F) so new strings are appended in the
# same order they were first used.
$ws.Range("C2:C4").Value = "Singapore Citizen"
$ws.Range("D2:D4").Value = 1
$ws.Range("E2:E4").Value = "EGDF21"
$ws.Range("F2:F4").Value = "IM2000"

# Column G no longer holds any data.
$ws.Range("G1:G4").ClearContents() | Out-Null

# Update the active selection to mirror the post-edit cursor position.
$ws.Range("D5").Select() | Out-Null

# Touch page setup (orientation) -- mirrors the pageSetup element added on save.
$ws.PageSetup.Orientation = 1
